# Data-driven update of Price (D) and Volume(1h) (E) columns, plus
# the three-row reshuffle (31-33) and two pair-swaps (45/46, 47/48)
# that resulted from the refreshed coin ranking on the source site.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue 'D2' '69.405.79'
Set-TextValue 'E2' '  +4.86%  '
Set-TextValue 'D3' '3.498.41'
Set-TextValue 'E3' '  +13.73%  '
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '588.98'
Set-TextValue 'E5' '  +2.59%  '
Set-TextValue 'D6' '184.94'
Set-TextValue 'E6' '  +8.59%  '
Set-TextValue 'D7' '3.499.59'
Set-TextValue 'E7' '  +13.90%  '
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.03%  '
Set-TextValue 'E9' '  +4.52%  '
Set-TextValue 'D10' '6.60'
Set-TextValue 'E10' '  +4.32%  '
Set-TextValue 'E11' '  +5.99%  '
Set-TextValue 'E12' '  +4.00%  '
Set-TextValue 'D13' '38.57'
Set-TextValue 'E13' '  +7.51%  '
Set-TextValue 'E14' '  +4.78%  '
Set-TextValue 'D15' '4.076.54'
Set-TextValue 'E15' '  +13.60%  '
Set-TextValue 'D16' '69.495.70'
Set-TextValue 'E16' '  +5.05%  '
Set-TextValue 'E17' '  +1.25%  '
Set-TextValue 'D18' '3.512.48'
Set-TextValue 'E18' '  +14.29%  '
Set-TextValue 'E19' '  +6.66%  '
Set-TextValue 'D20' '16.95'
Set-TextValue 'E20' '  +2.47%  '
Set-TextValue 'D21' '504.36'
Set-TextValue 'E21' '  +4.05%  '
Set-TextValue 'D22' '8.92'
Set-TextValue 'E22' '  +16.48%  '
Set-TextValue 'D23' '0.730'
Set-TextValue 'E23' '  +6.46%  '
Set-TextValue 'D24' '86.62'
Set-TextValue 'E24' '  +5.06%  '
Set-TextValue 'D25' '13.36'
Set-TextValue 'E25' '  +5.79%  '
Set-TextValue 'E26' '  +8.17%  '
Set-TextValue 'D27' '10.74'
Set-TextValue 'E27' '  +4.98%  '
Set-TextValue 'E28' '  +0.00%  '
Set-TextValue 'E29' '  +12.16%  '
Set-TextValue 'D30' '8.14'
Set-TextValue 'E30' '  +3.23%  '
Set-TextValue 'E34' '  +5.53%  '
Set-TextValue 'E35' '  +0.23%  '
Set-TextValue 'D36' '6.11'
Set-TextValue 'E36' '  +9.96%  '
Set-TextValue 'E37' '  +7.50%  '
Set-TextValue 'E38' '  +10.35%  '
Set-TextValue 'D39' '47.25'
Set-TextValue 'E39' '  +0.25%  '
Set-TextValue 'D40' '2.11'
Set-TextValue 'D41' '0.127'
Set-TextValue 'E41' '  +4.38%  '
Set-TextValue 'D42' '50.15'
Set-TextValue 'E42' '  +2.22%  '
Set-TextValue 'D43' '8.70'
Set-TextValue 'E43' '  +5.49%  '
Set-TextValue 'E44' '  +12.35%  '
Set-TextValue 'D49' '134.31'
Set-TextValue 'E49' '  -0.30%  '
Set-TextValue 'D51' '2.44'
Set-TextValue 'E51' '  +13.49%  '

# Row reshuffle: coin rank order changed for rows 31-33 and 45-48
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D31' '2.72'
Set-TextValue 'E31' '  +4.99%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D32' '0.0000107'
Set-TextValue 'E32' '  +19.58%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '30.70'
Set-TextValue 'E33' '  +10.63%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D45' '2.990.27'
Set-TextValue 'E45' '  +7.63%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D46' '407.14'
Set-TextValue 'E46' '  +11.93%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0363'
Set-TextValue 'E47' '  +5.78%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D48' '27.92'
Set-TextValue 'E48' '  +14.51%  '
